$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the whole working range to text format so numeric-looking values
# (hp, damage, ids, etc.) are preserved as text, matching the source schema.
$ws.Range("A4:Z11").NumberFormat = "@"

# Row 4
$ws.Range("A4").Value = "@"
$ws.Range("B4").Value = "@"
$ws.Range("C4").Value = "@"
$ws.Range("D4").Value = "string"
$ws.Range("E4").Value = "EnemyFamily"
$ws.Range("F4").Value = "uint"
$ws.Range("G4").Value = "uint"
$ws.Range("H4").Value = "float"
$ws.Range("I4").Value = "float"
$ws.Range("J4").Value = "AttackStyle"
$ws.Range("K4").Value = "float"
$ws.Range("L4").Value = "float"
$ws.Range("M4").Value = "float"
$ws.Range("N4").Value = "string"
$ws.Range("O4").Value = "string"
$ws.Range("P4").Value = "DamageType"
$ws.Range("Q4").Value = "DamageType"
$ws.Range("R4").Value = "string"
$ws.Range("S4").Value = "uint"
$ws.Range("T4").Value = "string"
$ws.Range("U4").Value = "uint"
$ws.Range("V4").Value = "string"
$ws.Range("W4").Value = "float"
$ws.Range("X4").Value = "string"
$ws.Range("Y4").Value = "string"
$ws.Range("Z4").Value = "string"

# Row 5
$ws.Range("A5").Value = "sector"
$ws.Range("B5").Value = "category"
$ws.Range("C5").Value = "serial"
$ws.Range("D5").Value = "name"
$ws.Range("E5").Value = "family"
$ws.Range("F5").Value = "hp"
$ws.Range("G5").Value = "damage"
$ws.Range("H5").Value = "moveSpeed"
$ws.Range("I5").Value = "radius"
$ws.Range("J5").Value = "attackStyle"
$ws.Range("K5").Value = "attackInterval"
$ws.Range("L5").Value = "projectileSpeed"
$ws.Range("M5").Value = "projectileLifetime"
$ws.Range("N5").Value = "projectileSprite"
$ws.Range("O5").Value = "impactSprite"
$ws.Range("P5").Value = "weakness"
$ws.Range("Q5").Value = "resistance"
$ws.Range("R5").Value = "lootTable"
$ws.Range("S5").Value = "sanityDamage"
$ws.Range("T5").Value = "combatNotes"
$ws.Range("U5").Value = "xp"
$ws.Range("V5").Value = "sprite"
$ws.Range("W5").Value = "spriteScale"
$ws.Range("X5").Value = "deathSprite"
$ws.Range("Y5").Value = "deathSfx"
$ws.Range("Z5").Value = "attackSfx"

# Row 6
$ws.Range("A6").Value = "40"
$ws.Range("B6").Value = "06"
$ws.Range("C6").Value = "0001"
$ws.Range("D6").Value = "裂界餍爬者"
$ws.Range("E6").Value = "SHAMBLER"
$ws.Range("F6").Value = "220"
$ws.Range("G6").Value = "24"
$ws.Range("H6").Value = "3.4"
$ws.Range("I6").Value = "12"
$ws.Range("J6").Value = "AUTO"
$ws.Range("K6").Value = "1.80"
$ws.Range("L6").Value = "14"
$ws.Range("M6").Value = "0.90"
$ws.Range("N6").Value = "fx/projectiles/spittle.png"
$ws.Range("O6").Value = "fx/impact/slime.png"
$ws.Range("P6").Value = "FIRE"
$ws.Range("Q6").Value = "VOID"
$ws.Range("R6").Value = "loot:ichor_minor"
$ws.Range("S6").Value = "6"
$ws.Range("T6").Value = "投掷腐质胆汁，落地后留下灼蚀雾。"
$ws.Range("U6").Value = "18"
$ws.Range("V6").Value = "ui/assets/topdown/top-down-shooter/characters/head/13.png"
$ws.Range("W6").Value = "0.9"
$ws.Range("X6").Value = "ui/assets/topdown/top-down-shooter/effects/explosion.png"
$ws.Range("Y6").Value = "ui/assets/topdown/top-down-shooter/sounds/explosion-2.wav"
$ws.Range("Z6").Value = "ui/assets/topdown/top-down-shooter/sounds/shoot-2.wav"

# Row 7
$ws.Range("A7").Value = "40"
$ws.Range("B7").Value = "06"
$ws.Range("C7").Value = "0002"
$ws.Range("D7").Value = "聆渊歌祭徒"
$ws.Range("E7").Value = "CULTIST"
$ws.Range("F7").Value = "260"
$ws.Range("G7").Value = "32"
$ws.Range("H7").Value = "3.1"
$ws.Range("I7").Value = "14"
$ws.Range("J7").Value = "BURST"
$ws.Range("K7").Value = "2.40"
$ws.Range("L7").Value = "22"
$ws.Range("M7").Value = "0.85"
$ws.Range("N7").Value = "fx/projectiles/choir_note.png"
$ws.Range("O7").Value = "fx/impact/chorus.png"
$ws.Range("P7").Value = "LIGHT"
$ws.Range("Q7").Value = "VOID"
$ws.Range("R7").Value = "loot:choir_cache"
$ws.Range("S7").Value = "9"
$ws.Range("T7").Value = "三连音符袭击，第三发附加畏惧层。"
$ws.Range("U7").Value = "26"
$ws.Range("V7").Value = "ui/assets/topdown/top-down-shooter/characters/head/7.png"
$ws.Range("W7").Value = "0.92"
$ws.Range("X7").Value = "ui/assets/topdown/top-down-shooter/effects/4.png"
$ws.Range("Y7").Value = "ui/assets/topdown/top-down-shooter/sounds/death.wav"
$ws.Range("Z7").Value = "ui/assets/topdown/top-down-shooter/sounds/shoot-3.wav"

# Row 8
$ws.Range("A8").Value = "40"
$ws.Range("B8").Value = "06"
$ws.Range("C8").Value = "0003"
$ws.Range("D8").Value = "渊喉呼嚎者"
$ws.Range("E8").Value = "ABERRATION"
$ws.Range("F8").Value = "340"
$ws.Range("G8").Value = "36"
$ws.Range("H8").Value = "4.0"
$ws.Range("I8").Value = "16"
$ws.Range("J8").Value = "MANUAL"
$ws.Range("K8").Value = "2.80"
$ws.Range("L8").Value = "0"
$ws.Range("M8").Value = "0.00"
$ws.Range("N8").Value = "fx/projectiles/howl_wave.png"
$ws.Range("O8").Value = "fx/impact/howl.png"
$ws.Range("P8").Value = "LIGHT"
$ws.Range("Q8").Value = "FROST"
$ws.Range("R8").Value = "loot:howler_pouch"
$ws.Range("S8").Value = "12"
$ws.Range("T8").Value = "扇形震荡波附加 3 秒理智流失。"
$ws.Range("U8").Value = "32"
$ws.Range("V8").Value = "ui/assets/topdown/top-down-shooter/characters/head/4.png"
$ws.Range("W8").Value = "0.95"
$ws.Range("X8").Value = "ui/assets/topdown/top-down-shooter/effects/3.png"
$ws.Range("Y8").Value = "ui/assets/topdown/top-down-shooter/sounds/explosion-3.wav"
$ws.Range("Z8").Value = "ui/assets/topdown/top-down-shooter/sounds/sword-2.wav"

# Row 9
$ws.Range("A9").Value = "40"
$ws.Range("B9").Value = "06"
$ws.Range("C9").Value = "0004"
$ws.Range("D9").Value = "虚壳哨兵"
$ws.Range("E9").Value = "CONSTRUCT"
$ws.Range("F9").Value = "420"
$ws.Range("G9").Value = "42"
$ws.Range("H9").Value = "2.6"
$ws.Range("I9").Value = "18"
$ws.Range("J9").Value = "BEAM"
$ws.Range("K9").Value = "1.45"
$ws.Range("L9").Value = "60"
$ws.Range("M9").Value = "0.50"
$ws.Range("N9").Value = "fx/projectiles/null_beam.png"
$ws.Range("O9").Value = "fx/impact/null_burn.png"
$ws.Range("P9").Value = "VOID"
$ws.Range("Q9").Value = "KINETIC"
$ws.Range("R9").Value = "loot:sentinel_cache"
$ws.Range("S9").Value = "10"
$ws.Range("T9").Value = "扫射光束前有 0.6 秒警示。"
$ws.Range("U9").Value = "38"
$ws.Range("V9").Value = "ui/assets/topdown/top-down-shooter/characters/turret/1.png"
$ws.Range("W9").Value = "1.05"
$ws.Range("X9").Value = "ui/assets/topdown/top-down-shooter/effects/5.png"
$ws.Range("Y9").Value = "ui/assets/topdown/top-down-shooter/sounds/explosion-1.wav"
$ws.Range("Z9").Value = "ui/assets/topdown/top-down-shooter/sounds/alert.wav"

# Row 10
$ws.Range("A10").Value = "40"
$ws.Range("B10").Value = "06"
$ws.Range("C10").Value = "0005"
$ws.Range("D10").Value = "堕港掘锚者"
$ws.Range("E10").Value = "CONSTRUCT"
$ws.Range("F10").Value = "520"
$ws.Range("G10").Value = "55"
$ws.Range("H10").Value = "2.2"
$ws.Range("I10").Value = "24"
$ws.Range("J10").Value = "MANUAL"
$ws.Range("K10").Value = "1.10"
$ws.Range("L10").Value = "0"
$ws.Range("M10").Value = "0.00"
$ws.Range("N10").Value = "fx/projectiles/dredger_slam.png"
$ws.Range("O10").Value = "fx/impact/dredger_slam.png"
$ws.Range("P10").Value = "FIRE"
$ws.Range("Q10").Value = "KINETIC"
$ws.Range("R10").Value = "loot:dredger_core"
$ws.Range("S10").Value = "14"
$ws.Range("T10").Value = "冲撞灯塔并引发地震波，需快速躲避。"
$ws.Range("U10").Value = "46"
$ws.Range("V10").Value = "ui/assets/topdown/top-down-shooter/characters/tank.png"
$ws.Range("W10").Value = "0.85"
$ws.Range("X10").Value = "ui/assets/topdown/top-down-shooter/effects/explosion.png"
$ws.Range("Y10").Value = "ui/assets/topdown/top-down-shooter/sounds/explosion-3.wav"
$ws.Range("Z10").Value = "ui/assets/topdown/top-down-shooter/sounds/flame-thrower.wav"

# Row 11
$ws.Range("A11").Value = "40"
$ws.Range("B11").Value = "06"
$ws.Range("C11").Value = "0006"
$ws.Range("D11").Value = "无数碎影"
$ws.Range("E11").Value = "ABERRATION"
$ws.Range("F11").Value = "160"
$ws.Range("G11").Value = "18"
$ws.Range("H11").Value = "4.8"
$ws.Range("I11").Value = "10"
$ws.Range("J11").Value = "BURST"
$ws.Range("K11").Value = "1.90"
$ws.Range("L11").Value = "26"
$ws.Range("M11").Value = "0.75"
$ws.Range("N11").Value = "fx/projectiles/fragment_dart.png"
$ws.Range("O11").Value = "fx/impact/fragment_spark.png"
$ws.Range("P11").Value = "LIGHT"
$ws.Range("Q11").Value = "VOID"
$ws.Range("R11").Value = "loot:fragment_cache"
$ws.Range("S11").Value = "8"
$ws.Range("T11").Value = "群猎碎片成群而行，形成交叉弹雨。"
$ws.Range("U11").Value = "20"
$ws.Range("V11").Value = "ui/assets/topdown/top-down-shooter/effects/1.png"
$ws.Range("W11").Value = "1.0"
$ws.Range("X11").Value = "ui/assets/topdown/top-down-shooter/effects/2.png"
$ws.Range("Y11").Value = "ui/assets/topdown/top-down-shooter/sounds/shoot-destroy.wav"
$ws.Range("Z11").Value = "ui/assets/topdown/top-down-shooter/sounds/shoot-1.wav"
